$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2629.25
$ws.Range("I100").Value = 2005
$ws.Range("J100").Value = 2837.3333
$ws.Range("K100").Value = 2005
$ws.Range("L100").Value = 2837.3333
$ws.Range("M100").Value = -1464
$ws.Range("N100").Value = -3919.3333
$ws.Range("H113").Value = 2852.75
$ws.Range("I113").Value = 2702.5
$ws.Range("J113").Value = 3003
$ws.Range("K113").Value = 2702.5
$ws.Range("L113").Value = 3003
$ws.Range("M113").Value = 551.5
$ws.Range("N113").Value = -9511
$ws.Range("H125").Value = 1632.6666
$ws.Range("J125").Value = 1449
$ws.Range("L125").Value = 13041
$ws.Range("N125").Value = -17961
$ws.Range("H127").Value = 1186.0667
$ws.Range("I127").Value = 458.8
$ws.Range("J127").Value = 1549.7
$ws.Range("K127").Value = 1376.4
$ws.Range("L127").Value = 4649.1
$ws.Range("M127").Value = 3583.6
$ws.Range("N127").Value = -14569.1
$ws.Range("H129").Value = 1000.9483
$ws.Range("J129").Value = 1054.151
$ws.Range("L129").Value = 3162.453
$ws.Range("N129").Value = -13162.453
$ws.Range("H132").Value = 708.95
$ws.Range("I132").Value = 658.27026
$ws.Range("J132").Value = 1334
$ws.Range("K132").Value = 1974.81078
$ws.Range("L132").Value = 4002
$ws.Range("M132").Value = 555.18922
$ws.Range("N132").Value = -9062
$ws.Range("H137").Value = 1856577.8
$ws.Range("I137").Value = 4691.6875
$ws.Range("K137").Value = 14075.0625
$ws.Range("M137").Value = -11525.0625

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H45").Value = 1331.6818
$ws.Range("I45").Value = 1234.9445
$ws.Range("J45").Value = 1767
$ws.Range("K45").Value = 1234.9445
$ws.Range("L45").Value = 1767
$ws.Range("M45").Value = -857.9445000000001
$ws.Range("N45").Value = -2521
$ws.Range("H52").Value = 53550
$ws.Range("J52").Value = 53550
$ws.Range("L52").Value = 53550
$ws.Range("N52").Value = -54186
$ws.Range("H61").Value = 6891.304
$ws.Range("I61").Value = 6595.6
$ws.Range("J61").Value = 8862.666999999999
$ws.Range("K61").Value = 6595.6
$ws.Range("L61").Value = 8862.666999999999
$ws.Range("M61").Value = -6383.6
$ws.Range("N61").Value = -9286.666999999999
$ws.Range("H110").Value = 1704.1666
$ws.Range("I110").Value = 1406.5
$ws.Range("J110").Value = 2299.5
$ws.Range("K110").Value = 1406.5
$ws.Range("L110").Value = 2299.5
$ws.Range("M110").Value = 638.5
$ws.Range("N110").Value = -6389.5
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H136").Value = 6891.304
$ws.Range("I136").Value = 6595.6
$ws.Range("J136").Value = 8862.666999999999
$ws.Range("K136").Value = 19786.8
$ws.Range("L136").Value = 26588.001
$ws.Range("M136").Value = -17236.8
$ws.Range("N136").Value = -31688.001
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = 0
$ws.Range("H141").Value = 66453.336
$ws.Range("J141").Value = 66453.336
$ws.Range("L141").Value = 66453.336
$ws.Range("N141").Value = -76813.336

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4323.9
$ws.Range("I105").Value = 4518
$ws.Range("K105").Value = 4518
$ws.Range("M105").Value = -2771
$ws.Range("H107").Value = 2614.8333
$ws.Range("I107").Value = 2664
$ws.Range("J107").Value = 2516.5
$ws.Range("K107").Value = 2664
$ws.Range("L107").Value = 2516.5
$ws.Range("M107").Value = -744
$ws.Range("N107").Value = -6356.5
$ws.Range("H134").Value = 7698.857
$ws.Range("I134").Value = 10277.5
$ws.Range("J134").Value = 6667.4
$ws.Range("K134").Value = 30832.5
$ws.Range("L134").Value = 20002.2
$ws.Range("M134").Value = -28297.5
$ws.Range("N134").Value = -25072.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 907623
$ws.Range("I31").Value = 7750
$ws.Range("J31").Value = 1426780.5
$ws.Range("K31").Value = 7750
$ws.Range("L31").Value = 1426780.5
$ws.Range("M31").Value = -7455
$ws.Range("N31").Value = -1427370.5
$ws.Range("H34").Value = 907623
$ws.Range("I34").Value = 7750
$ws.Range("J34").Value = 1426780.5
$ws.Range("K34").Value = 7750
$ws.Range("L34").Value = 1426780.5
$ws.Range("M34").Value = -7548
$ws.Range("N34").Value = -1427184.5
$ws.Range("H132").Value = 2958.8235
$ws.Range("I132").Value = 2917.8572
$ws.Range("J132").Value = 3150
$ws.Range("K132").Value = 8753.571599999999
$ws.Range("L132").Value = 9450
$ws.Range("M132").Value = -6223.571599999999
$ws.Range("N132").Value = -14510
$ws.Range("H134").Value = 2060.1277
$ws.Range("I134").Value = 1706.8889
$ws.Range("J134").Value = 3216.182
$ws.Range("K134").Value = 5120.6667
$ws.Range("L134").Value = 9648.545999999998
$ws.Range("M134").Value = -2585.6667
$ws.Range("N134").Value = -14718.546

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 5000
$ws.Range("J16").Value = 5000
$ws.Range("L16").Value = 15000
$ws.Range("N16").Value = -15346
$ws.Range("H68").Value = 2898.7585
$ws.Range("I68").Value = 964.5
$ws.Range("J68").Value = 4264.1177
$ws.Range("K68").Value = 2893.5
$ws.Range("L68").Value = 12792.3531
$ws.Range("M68").Value = -2082.5
$ws.Range("N68").Value = -14414.3531
$ws.Range("H71").Value = 2898.7585
$ws.Range("I71").Value = 964.5
$ws.Range("J71").Value = 4264.1177
$ws.Range("K71").Value = 8680.5
$ws.Range("L71").Value = 38377.0593
$ws.Range("M71").Value = -4624.5
$ws.Range("N71").Value = -46489.0593

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5756.915
$ws.Range("I70").Value = 5533.5356
$ws.Range("J70").Value = 6086.1055
$ws.Range("K70").Value = 5533.5356
$ws.Range("L70").Value = 6086.1055
$ws.Range("M70").Value = -5263.5356
$ws.Range("N70").Value = -6626.1055
$ws.Range("H73").Value = 5756.915
$ws.Range("I73").Value = 5533.5356
$ws.Range("J73").Value = 6086.1055
$ws.Range("K73").Value = 5533.5356
$ws.Range("L73").Value = 6086.1055
$ws.Range("M73").Value = -4597.5356
$ws.Range("N73").Value = -7958.1055

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1337.25
$ws.Range("I46").Value = 750
$ws.Range("J46").Value = 1924.5
$ws.Range("K46").Value = 750
$ws.Range("L46").Value = 1924.5
$ws.Range("M46").Value = -562
$ws.Range("N46").Value = -2300.5
$ws.Range("H94").Value = 27110
$ws.Range("J94").Value = 27110
$ws.Range("L94").Value = 27110
$ws.Range("N94").Value = -28462
$ws.Range("H136").Value = 4807
$ws.Range("I136").Value = 4629.3213
$ws.Range("K136").Value = 13887.9639
$ws.Range("M136").Value = -11337.9639

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 57499.75
$ws.Range("J121").Value = 57499.75
$ws.Range("L121").Value = 57499.75
$ws.Range("N121").Value = -60993.75
$ws.Range("H122").Value = 2800.8
$ws.Range("I122").Value = 2502
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 7506
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -5056
$ws.Range("N122").Value = -13900
$ws.Range("H126").Value = 1405.1765
$ws.Range("I126").Value = 1490.6666
$ws.Range("J126").Value = 1309
$ws.Range("K126").Value = 4471.9998
$ws.Range("L126").Value = 3927
$ws.Range("M126").Value = -2001.9998
$ws.Range("N126").Value = -8867
$ws.Range("H127").Value = 67166.664
$ws.Range("J127").Value = 67166.664
$ws.Range("L127").Value = 67166.664
$ws.Range("N127").Value = -77086.664
$ws.Range("H132").Value = 2515.2593
$ws.Range("I132").Value = 2500.2917
$ws.Range("K132").Value = 7500.875100000001
$ws.Range("M132").Value = -4970.875100000001
$ws.Range("H136").Value = 4745.42
$ws.Range("I136").Value = 1541.3914
$ws.Range("J136").Value = 7474.778
$ws.Range("K136").Value = 4624.174199999999
$ws.Range("L136").Value = 22424.334
$ws.Range("M136").Value = -2074.174199999999
$ws.Range("N136").Value = -27524.334
